$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.436.51'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '1.626.02'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.71'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.486'
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0618'
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("E10").Value = '  +3.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0830'
$ws.Range("E11").Value = '  +2.50%  '
$ws.Range("D12").Value = '1.850.98'
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '1.636.72'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '26.423.09'
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.64'
$ws.Range("E17").Value = '  +2.07%  '
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.48'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.35'
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.07'
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("E24").Value = '  -3.97%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.120'
$ws.Range("E27").Value = '  -2.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.22'
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0518'
$ws.Range("E30").Value = '  +4.74%  '
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("D36").Value = '1.157.82'
$ws.Range("E36").Value = '  +1.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0164'
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.806'
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.33'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.499'
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("E42").Value = '  +3.50%  '
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").Value = '1.762.37'
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.06'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.53'
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("E47").Value = '  +9.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.07'
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.410'
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  -0.15%  '
